$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters for E..T in order, matching the data dictionaries below
$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# A/D use shared strings (cluster names): 20=ECs, 21=FAPs, 22=sCs
$ECs = "ECs"
$FAPs = "FAPs"
$sCs = "sCs"

# Data for rows 2..10 : Sending(A), Target(D), then E..T values
$rows = @{
  2  = @{ A = $ECs;  D = $ECs;  V = @(3, 1, 88.72291666666666, 266.16875, 0.7675060578750151, 0.7675060578750152, 2, 0.6666666666666666, 1.321445333333333, 3.964336, 0.01021782062667047, 0.01021782062667047, 117.2424841888889, 1055.1823577, 0.007842239229249868, 0.00784223922924987) }
  3  = @{ A = $ECs;  D = $FAPs; V = @(3, 1, 88.72291666666666, 266.16875, 0.7675060578750151, 0.7675060578750152, 3, 1, 105.9632263333333, 317.889679, 0.819340166699254, 0.8193401666992541, 9401.366499703472, 84612.29849733124, 0.6288485414020022, 0.6288485414020024) }
  4  = @{ A = $ECs;  D = $sCs;  V = @(3, 1, 88.72291666666666, 266.16875, 0.7675060578750151, 0.7675060578750152, 3, 1, 22.04284166666666, 66.128525, 0.1704420126740755, 0.1704420126740755, 1955.705204288194, 17601.34683859375, 0.130815277243763, 0.1308152772437631) }
  5  = @{ A = $FAPs; D = $ECs;  V = @(3, 1, 17.91585, 53.74755, 0.1549827702197958, 0.1549827702197958, 2, 0.6666666666666666, 1.321445333333333, 3.964336, 0.01021782062667047, 0.01021782062667047, 23.6748163752, 213.0733473768, 0.001583586146330359, 0.00158358614633036) }
  6  = @{ A = $FAPs; D = $FAPs; V = @(3, 1, 17.91585, 53.74755, 0.1549827702197958, 0.1549827702197958, 3, 1, 105.9632263333333, 317.889679, 0.819340166699254, 0.8193401666992541, 1898.42126850405, 17085.79141653645, 0.1269836087873997, 0.1269836087873997) }
  7  = @{ A = $FAPs; D = $sCs;  V = @(3, 1, 17.91585, 53.74755, 0.1549827702197958, 0.1549827702197958, 3, 1, 22.04284166666666, 66.128525, 0.1704420126740755, 0.1704420126740755, 394.91624487375, 3554.24620386375, 0.02641557528606576, 0.02641557528606577) }
  8  = @{ A = $sCs;  D = $ECs;  V = @(3, 1, 8.960212333333333, 26.880637, 0.077511171905189, 0.07751117190518901, 2, 0.6666666666666666, 1.321445333333333, 3.964336, 0.01021782062667047, 0.01021782062667047, 11.84043077355911, 106.563876962032, 0.0007919952510902406, 0.0007919952510902409) }
  9  = @{ A = $sCs;  D = $FAPs; V = @(3, 1, 8.960212333333333, 26.880637, 0.077511171905189, 0.07751117190518901, 3, 1, 105.9632263333333, 317.889679, 0.819340166699254, 0.8193401666992541, 949.4530074717247, 8545.077067245524, 0.06350801650985209, 0.06350801650985211) }
  10 = @{ A = $sCs;  D = $sCs;  V = @(3, 1, 8.960212333333333, 26.880637, 0.077511171905189, 0.07751117190518901, 3, 1, 22.04284166666666, 66.128525, 0.1704420126740755, 0.1704420126740755, 197.5085417633805, 1777.576875870425, 0.01321116014424666, 0.01321116014424667) }
}

foreach ($r in 2..10) {
  $data = $rows[$r]
  $ws.Range("A$r").Value = $data.A
  $ws.Range("B$r").Value = "Tgfb1"
  $ws.Range("C$r").Value = "Sdc2"
  $ws.Range("D$r").Value = $data.D
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])$r").Value = $data.V[$i]
  }
}
